$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item(1)

# Update the SnippetId for the existing "getBcc" rows (19-20): the shared string
# "outlook-recipients-and-attendees-get-bcc-message-compose" is being renamed to
# "outlook-recipients-and-attendees-get-set-bcc-message-compose".
$ws.Range("D19").Value = "outlook-recipients-and-attendees-get-set-bcc-message-compose"
$ws.Range("D20").Value = "outlook-recipients-and-attendees-get-set-bcc-message-compose"

# Insert two new rows right after row 20 (before the old row 21) to host the new
# "setBcc" (Message Compose) snippet rows, shifting everything below down by 2.
$ws.Rows.Item(21).Resize(2).Insert()

# Populate the two new rows with the "set bcc" snippet metadata, mirroring the
# existing "set cc" rows pattern.
$ws.Range("A21").Value = "MessageCompose"
$ws.Range("B21").Value = "bcc"
$ws.Range("D21").Value = "outlook-recipients-and-attendees-get-set-bcc-message-compose"
$ws.Range("E21").Value = "setBcc"

$ws.Range("A22").Value = "Recipients"
$ws.Range("B22").Value = "setAsync"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "outlook-recipients-and-attendees-get-set-bcc-message-compose"
$ws.Range("E22").Value = "setBcc"

# Grow the "Snippets" table (and its autofilter) to cover the two new rows.
$tbl.Resize($ws.Range("A1:E158"))

# Match the recorded sheet view state after the edit.
$ws.Range("B1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E22").Select()
